# Insert a new data row at row 304 (pushes existing rows 304-409 down to 305-410)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(304).Insert()

# Populate the newly inserted row 304 with the new record's data
$ws.Range("A304").Value = 3
$ws.Range("B304").Value = "Femacal de La Calera"
$ws.Range("C304").Value = "Coquimbo"
$ws.Range("D304").Value = 44559
$ws.Range("E304").Value = 5
$ws.Range("F304").Value = 100114001
$ws.Range("G304").Value = "Papa"
$ws.Range("H304").Value = "Rosara"
$ws.Range("I304").Value = "1a (cosecha)"
$ws.Range("J304").Value = 600
$ws.Range("K304").Value = 8000
$ws.Range("L304").Value = 8500
$ws.Range("M304").Value = 8233
$ws.Range("N304").Value = "`$/saco 25 kilos"
$ws.Range("O304").Value = "Provincia de Quillota"
$ws.Range("P304").Value = 329
$ws.Range("Q304").Value = 25
$ws.Range("R304").Value = "Hortaliza"
